$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 15.493404
$ws.Range("H2").Value = 46.48021199999999
$ws.Range("I2").Value = 0.05356331879335558
$ws.Range("J2").Value = 0.05356331879335557
$ws.Range("O2").Value = 0.2013489143128838
$ws.Range("P2").Value = 0.2013489143128838
$ws.Range("Q2").Value = 19.486472532708
$ws.Range("R2").Value = 175.378252794372
$ws.Range("S2").Value = 0.01078491608603703
$ws.Range("T2").Value = 0.01078491608603703

# Row 3
$ws.Range("G3").Value = 15.493404
$ws.Range("H3").Value = 46.48021199999999
$ws.Range("I3").Value = 0.05356331879335558
$ws.Range("J3").Value = 0.05356331879335557
$ws.Range("M3").Value = 0.246708
$ws.Range("N3").Value = 0.740124
$ws.Range("O3").Value = 0.03949536580856015
$ws.Range("P3").Value = 0.03949536580856015
$ws.Range("Q3").Value = 3.822346714032
$ws.Range("R3").Value = 34.401120426288
$ws.Range("S3").Value = 0.002115502869664103
$ws.Range("T3").Value = 0.002115502869664103

# Row 4
$ws.Range("G4").Value = 15.493404
$ws.Range("H4").Value = 46.48021199999999
$ws.Range("I4").Value = 0.05356331879335558
$ws.Range("J4").Value = 0.05356331879335557
$ws.Range("M4").Value = 4.74207
$ws.Range("N4").Value = 14.22621
$ws.Range("O4").Value = 0.7591557198785561
$ws.Range("P4").Value = 0.759155719878556
$ws.Range("Q4").Value = 73.47080630628
$ws.Range("R4").Value = 661.23725675652
$ws.Range("S4").Value = 0.04066289983765445
$ws.Range("T4").Value = 0.04066289983765444

# Row 5
$ws.Range("I5").Value = 0.4524333485785276
$ws.Range("J5").Value = 0.4524333485785275
$ws.Range("O5").Value = 0.2013489143128838
$ws.Range("P5").Value = 0.2013489143128838
$ws.Range("S5").Value = 0.09109696353522902
$ws.Range("T5").Value = 0.091096963535229

# Row 6
$ws.Range("I6").Value = 0.4524333485785276
$ws.Range("J6").Value = 0.4524333485785275
$ws.Range("M6").Value = 0.246708
$ws.Range("N6").Value = 0.740124
$ws.Range("O6").Value = 0.03949536580856015
$ws.Range("P6").Value = 0.03949536580856015
$ws.Range("Q6").Value = 32.286220537032
$ws.Range("R6").Value = 290.575984833288
$ws.Range("S6").Value = 0.01786902060610075
$ws.Range("T6").Value = 0.01786902060610075

# Row 7
$ws.Range("I7").Value = 0.4524333485785276
$ws.Range("J7").Value = 0.4524333485785275
$ws.Range("M7").Value = 4.74207
$ws.Range("N7").Value = 14.22621
$ws.Range("O7").Value = 0.7591557198785561
$ws.Range("P7").Value = 0.759155719878556
$ws.Range("Q7").Value = 620.58594703878
$ws.Range("R7").Value = 5585.27352334902
$ws.Range("S7").Value = 0.3434673644371978
$ws.Range("T7").Value = 0.3434673644371977

# Row 8
$ws.Range("G8").Value = 66.835223
$ws.Range("H8").Value = 200.505669
$ws.Range("I8").Value = 0.2310606730563543
$ws.Range("J8").Value = 0.2310606730563542
$ws.Range("O8").Value = 0.2013489143128838
$ws.Range("P8").Value = 0.2013489143128838
$ws.Range("Q8").Value = 84.06046451812099
$ws.Range("R8").Value = 756.5441806630889
$ws.Range("S8").Value = 0.04652381566030113
$ws.Range("T8").Value = 0.04652381566030112

# Row 9
$ws.Range("G9").Value = 66.835223
$ws.Range("H9").Value = 200.505669
$ws.Range("I9").Value = 0.2310606730563543
$ws.Range("J9").Value = 0.2310606730563542
$ws.Range("M9").Value = 0.246708
$ws.Range("N9").Value = 0.740124
$ws.Range("O9").Value = 0.03949536580856015
$ws.Range("P9").Value = 0.03949536580856015
$ws.Range("Q9").Value = 16.488784195884
$ws.Range("R9").Value = 148.399057762956
$ws.Range("S9").Value = 0.009125825806332829
$ws.Range("T9").Value = 0.009125825806332829

# Row 10
$ws.Range("G10").Value = 66.835223
$ws.Range("H10").Value = 200.505669
$ws.Range("I10").Value = 0.2310606730563543
$ws.Range("J10").Value = 0.2310606730563542
$ws.Range("M10").Value = 4.74207
$ws.Range("N10").Value = 14.22621
$ws.Range("O10").Value = 0.7591557198785561
$ws.Range("P10").Value = 0.759155719878556
$ws.Range("Q10").Value = 316.93730593161
$ws.Range("R10").Value = 2852.43575338449
$ws.Range("S10").Value = 0.1754110315897203
$ws.Range("T10").Value = 0.1754110315897203

# Row 11
$ws.Range("G11").Value = 76.057215
$ws.Range("H11").Value = 228.171645
$ws.Range("I11").Value = 0.2629426595717627
$ws.Range("J11").Value = 0.2629426595717627
$ws.Range("O11").Value = 0.2013489143128838
$ws.Range("P11").Value = 0.2013489143128838
$ws.Range("Q11").Value = 95.659212850305
$ws.Range("R11").Value = 860.9329156527451
$ws.Range("S11").Value = 0.05294321903131662
$ws.Range("T11").Value = 0.05294321903131662

# Row 12
$ws.Range("G12").Value = 76.057215
$ws.Range("H12").Value = 228.171645
$ws.Range("I12").Value = 0.2629426595717627
$ws.Range("J12").Value = 0.2629426595717627
$ws.Range("M12").Value = 0.246708
$ws.Range("N12").Value = 0.740124
$ws.Range("O12").Value = 0.03949536580856015
$ws.Range("P12").Value = 0.03949536580856015
$ws.Range("Q12").Value = 18.76392339822
$ws.Range("R12").Value = 168.87531058398
$ws.Range("S12").Value = 0.01038501652646247
$ws.Range("T12").Value = 0.01038501652646247

# Row 13
$ws.Range("G13").Value = 76.057215
$ws.Range("H13").Value = 228.171645
$ws.Range("I13").Value = 0.2629426595717627
$ws.Range("J13").Value = 0.2629426595717627
$ws.Range("M13").Value = 4.74207
$ws.Range("N13").Value = 14.22621
$ws.Range("O13").Value = 0.7591557198785561
$ws.Range("P13").Value = 0.759155719878556
$ws.Range("Q13").Value = 360.66863753505
$ws.Range("R13").Value = 3246.01773781545
$ws.Range("S13").Value = 0.1996144240139836
$ws.Range("T13").Value = 0.1996144240139836
